$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ciudades")

# Swap the Huesca/Huelva rows (row 53 and row 54): the name and the
# full set of numeric stats swap places between the two rows.
$nameA = $ws.Range("A53").Value2
$nameB = $ws.Range("A54").Value2
$ws.Range("A53").Value2 = $nameB
$ws.Range("A54").Value2 = $nameA

$bA = $ws.Range("B53").Value2
$cA = $ws.Range("C53").Value2
$dA = $ws.Range("D53").Value2
$eA = $ws.Range("E53").Value2

$bB = $ws.Range("B54").Value2
$cB = $ws.Range("C54").Value2
$dB = $ws.Range("D54").Value2
$eB = $ws.Range("E54").Value2

$ws.Range("B53").Value2 = $bB
$ws.Range("C53").Value2 = $cB
$ws.Range("D53").Value2 = $dB
$ws.Range("E53").Value2 = $eB

$ws.Range("B54").Value2 = $bA
$ws.Range("C54").Value2 = $cA
$ws.Range("D54").Value2 = $dA
$ws.Range("E54").Value2 = $eA

# Update the "last updated" timestamp in cell A1.
$ws.Range("A1").Value2 = "Datos actualizados a 22 de Marzo de 2020 a las 09:16"
